# Added exception handling to handle hangs when the vpn isn't working well,
# plus a timeout on the stream reader. New "Random" method rows captured
# during this run are appended to the IBB data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - new sample captured from the Random method
$ws.Range("A3").Value = 42600.829224537039
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("B3").Value = "Random"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 77
$ws.Range("I3").Value = 23
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 27
$ws.Range("M3").Value = 73

# Row 4 - new sample captured from the Random method
$ws.Range("A4").Value = 42600.882037037038
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = "Random"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 85
$ws.Range("I4").Value = 15
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 45
$ws.Range("M4").Value = 55

# Column A needs to widen slightly to keep fitting the date/time values
$ws.Columns("A:A").ColumnWidth = 14
